$d = $word.ActiveDocument

# Step 1: the run holding the placeholder gains a trailing space, so the
# word "anos" (now implied inside max_pena's own rendered value) can be
# dropped from the following run without losing the separating space.
#   "{{ max_pena }}"  ->  "{{ max_pena }} "
$d.Content.Find.Execute(
    "{{ max_pena }}",   # FindText
    $true,              # MatchCase
    $false,             # MatchWholeWord
    $false,             # MatchWildcards
    $false,             # MatchSoundsLike
    $false,             # MatchAllWordForms
    $true,              # Forward
    1,                  # Wrap (wdFindContinue)
    $false,             # Format
    "{{ max_pena }} ",  # ReplaceWith
    2)                  # Replace (wdReplaceAll)

# Step 2: drop the leading " anos " from the run that follows, leaving
# "de reclusão, consoante art. 109 do Código Penal."
$d.Content.Find.Execute(
    " anos de reclusão, consoante",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "de reclusão, consoante",
    2)
